# Fix Barclay XLS format
#
# Renames a handful of column-header / repeated-value labels in the
# "Barclaycard Umsaetze" sheet (shortening the German merchant-detail
# wording), and moves the saved cursor/selection back to the top of the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 13) relabeling ------------------------------------
$ws.Range("L13").Value2 = "Karteninhaber"   # was "Name des Karteninhabers"
$ws.Range("O13").Value2 = "Details"         # was "Händlerdetails"

# --- "Händler<N>" -> "Detail<N>" for every data row (14-20) ------------
$ws.Range("O14").Value2 = "Detail0"         # was "Händler0"
$ws.Range("O15").Value2 = "Detail1"         # was "Händler1"
$ws.Range("O16").Value2 = "Detail2"         # was "Händler2"
$ws.Range("O17").Value2 = "Detail3"         # was "Händler3"
$ws.Range("O18").Value2 = "Detail4"         # was "Händler4"
$ws.Range("O19").Value2 = "Detail5"         # was "Händler5"
$ws.Range("O20").Value2 = "Detail6"         # was "Händler6"

# --- Reset the saved selection/active cell from K14 to E5 --------------
$ws.Range("E5").Select() | Out-Null
